# Insert a new data row at row 483 (pushing the existing rows 483..586
# down to 484..587) and populate it with the new record.
#
# Columns (row 1 headers):
#   A Mercado ID            B Mercado                 C Región
#   D Fecha                 E Codreg                  F Categoría ID
#   G Categoría              H Variedad                 I Calidad
#   J Volumen                K Precio mínimo            L Precio máximo
#   M Precio promedio ponderado   N Unidad de comercialización
#   O Origen                 P Precio $/Kg              Q Kg o Unidades
#   R Clasificación

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 483:586 down to 484:587, leaving a blank row 483.
$ws.Rows("483:483").Insert()

$row = 483

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"
$ws.Cells.Item($row, 4).Value = 44641
$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = 100112020
$ws.Cells.Item($row, 7).Value = "Tomate"
$ws.Cells.Item($row, 8).Value = "Larga vida"
$ws.Cells.Item($row, 9).Value = "Extra"
$ws.Cells.Item($row, 10).Value = 300
$ws.Cells.Item($row, 11).Value = 25000
$ws.Cells.Item($row, 12).Value = 25000
$ws.Cells.Item($row, 13).Value = 25000
$ws.Cells.Item($row, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($row, 16).Value = 1389
$ws.Cells.Item($row, 17).Value = 18
$ws.Cells.Item($row, 18).Value = "Hortaliza"
